# Auto-generated edit script: updates computed profit-margin columns (H:N)
# in the Hyperion_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to refreshed market-board snapshot values.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 262.33334
$ws.Range("I2").Value = 348.33334
$ws.Range("J2").Value = 90.333336
$ws.Range("K2").Value = 348.33334
$ws.Range("L2").Value = 90.333336
$ws.Range("M2").Value = -235.33334
$ws.Range("N2").Value = -316.333336
$ws.Range("H17").Value = 1218.5312
$ws.Range("I17").Value = 634.5
$ws.Range("J17").Value = 1301.9642
$ws.Range("K17").Value = 1903.5
$ws.Range("L17").Value = 3905.8926
$ws.Range("M17").Value = -1735.5
$ws.Range("N17").Value = -4241.892599999999
$ws.Range("H19").Value = 5599.143
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 5599.143
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 5599.143
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -5949.143
$ws.Range("H32").Value = 3299.4443
$ws.Range("I32").Value = 1450
$ws.Range("J32").Value = 3699.3242
$ws.Range("K32").Value = 1450
$ws.Range("L32").Value = 3699.3242
$ws.Range("M32").Value = -1124
$ws.Range("N32").Value = -4351.3242
$ws.Range("H40").Value = 4216.407
$ws.Range("I40").Value = 1489
$ws.Range("J40").Value = 4434.6
$ws.Range("K40").Value = 1489
$ws.Range("L40").Value = 4434.6
$ws.Range("M40").Value = -1314
$ws.Range("N40").Value = -4784.6
$ws.Range("H55").Value = 244.90909
$ws.Range("I55").Value = 165.2
$ws.Range("J55").Value = 311.33334
$ws.Range("K55").Value = 165.2
$ws.Range("L55").Value = 311.33334
$ws.Range("M55").Value = 48.80000000000001
$ws.Range("N55").Value = -739.33334
$ws.Range("H62").Value = 8689.546
$ws.Range("I62").Value = 8420.556
$ws.Range("K62").Value = 8420.556
$ws.Range("M62").Value = -7796.556
$ws.Range("H65").Value = 8689.546
$ws.Range("I65").Value = 8420.556
$ws.Range("K65").Value = 42102.78
$ws.Range("M65").Value = -38982.78
$ws.Range("H74").Value = 7532.1333
$ws.Range("I74").Value = 7000
$ws.Range("J74").Value = 7570.143
$ws.Range("K74").Value = 7000
$ws.Range("L74").Value = 7570.143
$ws.Range("M74").Value = -6064
$ws.Range("N74").Value = -9442.143
$ws.Range("H77").Value = 7532.1333
$ws.Range("I77").Value = 7000
$ws.Range("J77").Value = 7570.143
$ws.Range("K77").Value = 35000
$ws.Range("L77").Value = 37850.715
$ws.Range("M77").Value = -30320
$ws.Range("N77").Value = -47210.715
$ws.Range("H86").Value = 2042.3448
$ws.Range("I86").Value = 2081.6
$ws.Range("J86").Value = 1955.1111
$ws.Range("K86").Value = 2081.6
$ws.Range("L86").Value = 1955.1111
$ws.Range("M86").Value = -958.5999999999999
$ws.Range("N86").Value = -4201.1111
$ws.Range("H89").Value = 2042.3448
$ws.Range("I89").Value = 2081.6
$ws.Range("J89").Value = 1955.1111
$ws.Range("K89").Value = 10408
$ws.Range("L89").Value = 9775.5555
$ws.Range("M89").Value = -4792
$ws.Range("N89").Value = -21007.5555
$ws.Range("H97").Value = 1496.4
$ws.Range("J97").Value = 1370.5
$ws.Range("L97").Value = 4111.5
$ws.Range("N97").Value = -5103.5
$ws.Range("H98").Value = 1846.5264
$ws.Range("I98").Value = 1805.6666
$ws.Range("K98").Value = 1805.6666
$ws.Range("M98").Value = -307.6666
$ws.Range("H122").Value = 1846.5264
$ws.Range("I122").Value = 1805.6666
$ws.Range("K122").Value = 5416.9998
$ws.Range("M122").Value = -2966.9998
$ws.Range("H132").Value = 22730470
$ws.Range("I132").Value = 23812754
$ws.Range("K132").Value = 71438262
$ws.Range("M132").Value = -71435732
$ws.Range("H133").Value = 69833.2
$ws.Range("H136").Value = 126292
$ws.Range("J136").Value = 182589
$ws.Range("L136").Value = 182589
$ws.Range("N136").Value = -192789

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 25342.25
$ws.Range("J24").Value = 25342.25
$ws.Range("L24").Value = 25342.25
$ws.Range("N24").Value = -26090.25
$ws.Range("H32").Value = 3256.4893
$ws.Range("I32").Value = 1889.7866
$ws.Range("J32").Value = 8651.368
$ws.Range("K32").Value = 1889.7866
$ws.Range("L32").Value = 8651.368
$ws.Range("M32").Value = -1602.7866
$ws.Range("N32").Value = -9225.368
$ws.Range("H61").Value = 2674.1853
$ws.Range("I61").Value = 2422.348
$ws.Range("K61").Value = 2422.348
$ws.Range("M61").Value = -2210.348
$ws.Range("H75").Value = 15999
$ws.Range("I75").Value = 15999
$ws.Range("K75").Value = 15999
$ws.Range("M75").Value = -15125
$ws.Range("H78").Value = 15999
$ws.Range("I78").Value = 15999
$ws.Range("K78").Value = 47997
$ws.Range("M78").Value = -43629
$ws.Range("H97").Value = 1205576.5
$ws.Range("I97").Value = 1707565
$ws.Range("K97").Value = 1707565
$ws.Range("M97").Value = -1707069
$ws.Range("H100").Value = 25342.25
$ws.Range("J100").Value = 25342.25
$ws.Range("L100").Value = 25342.25
$ws.Range("N100").Value = -27506.25
$ws.Range("H110").Value = 6951942.5
$ws.Range("I110").Value = 9265924
$ws.Range("J110").Value = 9999
$ws.Range("K110").Value = 9265924
$ws.Range("L110").Value = 9999
$ws.Range("M110").Value = -9263879
$ws.Range("N110").Value = -14089
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
$ws.Range("H127").Value = 96666.664
$ws.Range("J127").Value = 96666.664
$ws.Range("L127").Value = 96666.664
$ws.Range("N127").Value = -106586.664
$ws.Range("H132").Value = 2381.95
$ws.Range("I132").Value = 1603.125
$ws.Range("J132").Value = 5497.25
$ws.Range("K132").Value = 4809.375
$ws.Range("L132").Value = 16491.75
$ws.Range("M132").Value = -2279.375
$ws.Range("N132").Value = -21551.75
$ws.Range("H136").Value = 2674.1853
$ws.Range("I136").Value = 2422.348
$ws.Range("K136").Value = 7267.044
$ws.Range("M136").Value = -4717.044

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3968426.5
$ws.Range("I22").Value = 3968426.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3968426.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3968253.5
$ws.Range("N22").ClearContents()
$ws.Range("H33").Value = 9500
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 9500
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 9500
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -10172
$ws.Range("H86").Value = 5006785.5
$ws.Range("I86").Value = 12513927
$ws.Range("K86").Value = 12513927
$ws.Range("M86").Value = -12512804
$ws.Range("H89").Value = 5006785.5
$ws.Range("I89").Value = 12513927
$ws.Range("K89").Value = 62569635
$ws.Range("M89").Value = -62564019
$ws.Range("H94").Value = 2939285.5
$ws.Range("I94").Value = 4329774.5
$ws.Range("J94").Value = 19258.8
$ws.Range("K94").Value = 4329774.5
$ws.Range("L94").Value = 19258.8
$ws.Range("M94").Value = -4329323.5
$ws.Range("N94").Value = -20160.8
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H127").Value = 44000
$ws.Range("J127").Value = 44000
$ws.Range("L127").Value = 44000
$ws.Range("N127").Value = -53920
$ws.Range("H134").Value = 7199.5386
$ws.Range("I134").Value = 3149.875
$ws.Range("J134").Value = 13679
$ws.Range("K134").Value = 9449.625
$ws.Range("L134").Value = 41037
$ws.Range("M134").Value = -6914.625
$ws.Range("N134").Value = -46107

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20872.021
$ws.Range("I31").Value = 1405.4333
$ws.Range("J31").Value = 53316.332
$ws.Range("K31").Value = 1405.4333
$ws.Range("L31").Value = 53316.332
$ws.Range("M31").Value = -1110.4333
$ws.Range("N31").Value = -53906.332
$ws.Range("H34").Value = 20872.021
$ws.Range("I34").Value = 1405.4333
$ws.Range("J34").Value = 53316.332
$ws.Range("K34").Value = 1405.4333
$ws.Range("L34").Value = 53316.332
$ws.Range("M34").Value = -1203.4333
$ws.Range("N34").Value = -53720.332
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250
$ws.Range("H51").Value = 9999
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 2996.6592
$ws.Range("I58").Value = 3212.7812
$ws.Range("J58").Value = 2420.3333
$ws.Range("K58").Value = 3212.7812
$ws.Range("L58").Value = 2420.3333
$ws.Range("M58").Value = -3009.7812
$ws.Range("N58").Value = -2826.3333
$ws.Range("H59").Value = 39999.5
$ws.Range("J59").Value = 39999.5
$ws.Range("L59").Value = 39999.5
$ws.Range("N59").Value = -42289.5
$ws.Range("H61").Value = 9999
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 1999.8334
$ws.Range("I62").Value = 1499.75
$ws.Range("K62").Value = 1499.75
$ws.Range("M62").Value = -875.75
$ws.Range("H65").Value = 1999.8334
$ws.Range("I65").Value = 1499.75
$ws.Range("K65").Value = 7498.75
$ws.Range("M65").Value = -4378.75
$ws.Range("H88").Value = 49495.5
$ws.Range("I88").Value = 25999
$ws.Range("K88").Value = 25999
$ws.Range("M88").Value = -25593
$ws.Range("H91").Value = 49495.5
$ws.Range("I91").Value = 25999
$ws.Range("K91").Value = 25999
$ws.Range("M91").Value = -24595
$ws.Range("H94").Value = 1012.9474
$ws.Range("J94").Value = 1163.3077
$ws.Range("L94").Value = 1163.3077
$ws.Range("N94").Value = -2065.3077
$ws.Range("H97").Value = 28498
$ws.Range("J97").Value = 28498
$ws.Range("L97").Value = 28498
$ws.Range("N97").Value = -30480
$ws.Range("H99").Value = 4299.4
$ws.Range("J99").Value = 3999.25
$ws.Range("L99").Value = 3999.25
$ws.Range("N99").Value = -6995.25
$ws.Range("H105").Value = 1537.4736
$ws.Range("I105").Value = 1302.1333
$ws.Range("J105").Value = 2420
$ws.Range("K105").Value = 1302.1333
$ws.Range("L105").Value = 2420
$ws.Range("M105").Value = 444.8667
$ws.Range("N105").Value = -5914
$ws.Range("H109").Value = 36830.168
$ws.Range("J109").Value = 41196.2
$ws.Range("L109").Value = 41196.2
$ws.Range("N109").Value = -43276.2
$ws.Range("H126").Value = 4299.4
$ws.Range("J126").Value = 3999.25
$ws.Range("L126").Value = 11997.75
$ws.Range("N126").Value = -16937.75
$ws.Range("H132").Value = 46604.65
$ws.Range("I132").Value = 1718.2354
$ws.Range("K132").Value = 5154.706200000001
$ws.Range("M132").Value = -2624.706200000001
$ws.Range("H134").Value = 3205.1667
$ws.Range("I134").Value = 2606.6667
$ws.Range("J134").Value = 4402.1665
$ws.Range("K134").Value = 7820.000100000001
$ws.Range("L134").Value = 13206.4995
$ws.Range("M134").Value = -5285.000100000001
$ws.Range("N134").Value = -18276.4995
$ws.Range("H136").Value = 2996.6592
$ws.Range("I136").Value = 3212.7812
$ws.Range("J136").Value = 2420.3333
$ws.Range("K136").Value = 9638.3436
$ws.Range("L136").Value = 7260.999899999999
$ws.Range("M136").Value = -7088.3436
$ws.Range("N136").Value = -12360.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9264873
$ws.Range("I56").Value = 9264873
$ws.Range("K56").Value = 9264873
$ws.Range("M56").Value = -9264343
$ws.Range("H75").Value = 2836
$ws.Range("J75").Value = 3470
$ws.Range("L75").Value = 10410
$ws.Range("N75").Value = -12406
$ws.Range("H78").Value = 2836
$ws.Range("J78").Value = 3470
$ws.Range("L78").Value = 31230
$ws.Range("N78").Value = -41214
$ws.Range("H131").Value = 13891003
$ws.Range("I131").Value = 9259952
$ws.Range("J131").Value = 15875739
$ws.Range("K131").Value = 27779856
$ws.Range("L131").Value = 47627217
$ws.Range("M131").Value = -27774816
$ws.Range("N131").Value = -47637297

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4352479
$ws.Range("I70").Value = 5559810.5
$ws.Range("K70").Value = 5559810.5
$ws.Range("M70").Value = -5559540.5
$ws.Range("H73").Value = 4352479
$ws.Range("I73").Value = 5559810.5
$ws.Range("K73").Value = 5559810.5
$ws.Range("M73").Value = -5558874.5
$ws.Range("H95").Value = 49997.5
$ws.Range("J95").Value = 49997.5
$ws.Range("L95").Value = 49997.5
$ws.Range("N95").Value = -55489.5
$ws.Range("H97").Value = 993187.4399999999
$ws.Range("I97").Value = 1191584.4
$ws.Range("J97").Value = 1202.75
$ws.Range("K97").Value = 1191584.4
$ws.Range("L97").Value = 1202.75
$ws.Range("M97").Value = -1191088.4
$ws.Range("N97").Value = -2194.75
$ws.Range("H98").Value = 40000
$ws.Range("J98").Value = 40000
$ws.Range("L98").Value = 40000
$ws.Range("N98").Value = -45990
$ws.Range("H122").Value = 427623.9
$ws.Range("I122").Value = 686467.6
$ws.Range("K122").Value = 2059402.8
$ws.Range("M122").Value = -2056952.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4229.45
$ws.Range("I7").Value = 2816.6155
$ws.Range("K7").Value = 2816.6155
$ws.Range("M7").Value = -2704.6155
$ws.Range("H16").Value = 972.71875
$ws.Range("I16").Value = 771.08
$ws.Range("J16").Value = 1692.8572
$ws.Range("K16").Value = 771.08
$ws.Range("L16").Value = 1692.8572
$ws.Range("M16").Value = -601.08
$ws.Range("N16").Value = -2032.8572
$ws.Range("H40").Value = 10629.5
$ws.Range("I40").Value = 8860
$ws.Range("K40").Value = 8860
$ws.Range("M40").Value = -8724
$ws.Range("H46").Value = 5305.8887
$ws.Range("I46").Value = 5092.1816
$ws.Range("J46").Value = 5452.8125
$ws.Range("K46").Value = 5092.1816
$ws.Range("L46").Value = 5452.8125
$ws.Range("M46").Value = -4904.1816
$ws.Range("N46").Value = -5828.8125
$ws.Range("H68").Value = 3010.4546
$ws.Range("I68").Value = 2639.5
$ws.Range("K68").Value = 2639.5
$ws.Range("M68").Value = -1890.5
$ws.Range("H71").Value = 3010.4546
$ws.Range("I71").Value = 2639.5
$ws.Range("K71").Value = 13197.5
$ws.Range("M71").Value = -9453.5
$ws.Range("H96").Value = 28000
$ws.Range("J96").Value = 28000
$ws.Range("L96").Value = 28000
$ws.Range("N96").Value = -33492
$ws.Range("H101").Value = 15177.5
$ws.Range("J101").Value = 15177.5
$ws.Range("L101").Value = 15177.5
$ws.Range("N101").Value = -21667.5
$ws.Range("H122").Value = 5323.1333
$ws.Range("I122").Value = 3166.4
$ws.Range("J122").Value = 9636.6
$ws.Range("K122").Value = 9499.200000000001
$ws.Range("L122").Value = 28909.8
$ws.Range("M122").Value = -7049.200000000001
$ws.Range("N122").Value = -33809.8
$ws.Range("H126").Value = 4229.45
$ws.Range("I126").Value = 2816.6155
$ws.Range("K126").Value = 8449.8465
$ws.Range("M126").Value = -5979.8465
$ws.Range("H132").Value = 4287.52
$ws.Range("I132").Value = 3728.111
$ws.Range("K132").Value = 11184.333
$ws.Range("M132").Value = -8654.332999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 5517.5
$ws.Range("I51").Value = 5517.5
$ws.Range("K51").Value = 5517.5
$ws.Range("M51").Value = -5007.5
$ws.Range("H52").Value = 8562.5
$ws.Range("I52").Value = 5000
$ws.Range("K52").Value = 5000
$ws.Range("M52").Value = -4774
$ws.Range("H107").Value = 45456140
$ws.Range("I107").Value = 47620640
$ws.Range("K107").Value = 142861920
$ws.Range("M107").Value = -142860000
$ws.Range("H108").Value = 69994
$ws.Range("J108").Value = 69994
$ws.Range("L108").Value = 69994
$ws.Range("N108").Value = -77674
$ws.Range("H109").Value = 35994.5
$ws.Range("J109").Value = 35994.5
$ws.Range("L109").Value = 35994.5
$ws.Range("N109").Value = -38768.5
$ws.Range("H123").Value = 60476
$ws.Range("J123").Value = 60476
$ws.Range("L123").Value = 60476
$ws.Range("N123").Value = -70276
